# Adds the new observation row (row 5) to the "Artfynd" sheet, matching
# the source system export for Id 112330395 (Goliatmusseron / Tricholoma
# matsutake found at Skrakaholmberget, Ang on 2023-09-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

# --- numeric columns -------------------------------------------------
$ws.Cells.Item($row, 1).Value  = 112330395    # A  Id
$ws.Cells.Item($row, 2).Value  = 88166        # B  Taxonsorteringsordning
$ws.Cells.Item($row, 5).Value  = 6276         # E  TaxonId
$ws.Cells.Item($row, 17).Value = 560633       # Q  Ost
$ws.Cells.Item($row, 18).Value = 7108638      # R  Nord
$ws.Cells.Item($row, 19).Value = 10           # S  Noggrannhet

# --- text columns ------------------------------------------------------
$ws.Cells.Item($row, 3).Value  = "Ovaliderad"                     # C  Valideringsstatus
$ws.Cells.Item($row, 4).Value  = "VU"                              # D  Rödlistade
$ws.Cells.Item($row, 6).Value  = "Goliatmusseron"                  # F  Artnamn
$ws.Cells.Item($row, 7).Value  = "Tricholoma matsutake"             # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = "(S.Ito & S.Imai) Singer"         # H  Auktor
$ws.Cells.Item($row, 16).Value = "Skrakaholmberget, Ång"           # P  Lokalnamn
$ws.Cells.Item($row, 20).Value = "Jämtland"                        # T  Län
$ws.Cells.Item($row, 21).Value = "Strömsund"                       # U  Kommun
$ws.Cells.Item($row, 22).Value = "Ångermanland"                    # V  Provins
$ws.Cells.Item($row, 23).Value = "Tåsjö"                           # W  Församling
$ws.Cells.Item($row, 49).Value = "Susanne Wiik"                    # AW Rapportör
$ws.Cells.Item($row, 50).Value = "Susanne Wiik"                    # AX Observatörer

# Dates are stored as plain text (matching the other rows on this sheet),
# not as native Excel date serials, so the cells are pre-formatted as Text
# before the value is entered to keep them as strings (not auto-converted
# date serials) on save.
$ws.Cells.Item($row, 25).NumberFormat = "@"   # Y  Startdatum
$ws.Cells.Item($row, 25).Value = "2023-09-12"
$ws.Cells.Item($row, 27).NumberFormat = "@"   # AA Slutdatum
$ws.Cells.Item($row, 27).Value = "2023-09-12"

# --- boolean columns ----------------------------------------------------
$ws.Cells.Item($row, 30).Value = $false   # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false   # AE Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false   # AG Ospontan
